$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "244.86"
Set-TextValue "D4"  "5.215"
Set-TextValue "D5"  "0.05784"
Set-TextValue "D6"  "6.498"
Set-TextValue "D7"  "3.130"
Set-TextValue "D8"  "0.8167"
Set-TextValue "D9"  "0.8555"
Set-TextValue "D10" "0.1362"
Set-TextValue "D11" "0.06968"
Set-TextValue "D12" "0.03177"
Set-TextValue "D13" "0.02873"
Set-TextValue "D14" "0.09377"
Set-TextValue "D15" "3.748"
Set-TextValue "D16" "0.001522"
Set-TextValue "D17" "0.04699"
Set-TextValue "D18" "0.0005983"
Set-TextValue "D19" "0.006278"
Set-TextValue "D20" "0.001237"
Set-TextValue "D21" "0.004527"
Set-TextValue "D22" "0.00006910"
Set-TextValue "D23" "3.501"
Set-TextValue "D25" "0.3175"
Set-TextValue "D26" "0.1338"
Set-TextValue "D27" "0.1327"
Set-TextValue "D28" "0.0002330"
Set-TextValue "D40" "0.03655"

Set-TextValue "D41" "0.006309"
$ws.Range("E41").Value = "40KickTokenKICK"

Set-TextValue "D42" "0.1053"

Set-TextValue "D43" "0.002684"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

Set-TextValue "D44" "0.008322"
Set-TextValue "D45" "0.00005271"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

Set-TextValue "D48" "0.002349"
$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"
